{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== 11) {\n  throw new Error(\"Unexpected paragraph count: \" + paragraphs.items.length);\n}\n\nparagraphs.items[0].insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 24.05.25\\u000brStar-Math: Small LLMs Can Master Math Reasoning with Self-Evolved Deep Thinking\", Word.InsertLocation.replace);\nparagraphs.items[1].insertText(\"\u05db\u05de\u05d4 \u05d9\u05de\u05d9\u05dd \u05dc\u05d0 \u05e1\u05e7\u05e8\u05ea\u05d9 \u05de\u05d0\u05de\u05e8 \u05d0\u05d1\u05dc \u05d1\u05d9\u05d5\u05dd \u05d4\u05d5\u05dc\u05d3\u05ea\u05d9 \u05dc\u05d0 \u05d9\u05db\u05d5\u05dc\u05ea\u05d9 \u05dc\u05d0 \u05dc\u05db\u05ea\u05d5\u05d1 \u05e1\u05e7\u05d9\u05e8\u05d4 \u05dc\u05de\u05e8\u05d5\u05ea \u05d4\u05e2\u05d5\u05de\u05e1 \u05d4\u05de\u05d8\u05d5\u05e8\u05e3. \u05d4\u05d9\u05d5\u05dd \u05d0\u05e1\u05e7\u05d5\u05e8 \u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05d9\u05e6\u05d0 \u05dc\u05e4\u05e0\u05d9 4 \u05d7\u05d5\u05d3\u05e9\u05d9\u05dd \u05d5\u05d4\u05d5\u05d0 \u05de\u05e9\u05dc\u05d1 \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea MCTS \u05e9\u05d6\u05d4 \u05e7\u05d9\u05e6\u05d5\u05e8 \u05e9\u05dc Monte Carlo Tree Search. \u05e8\u05d5\u05d1\u05db\u05dd \u05db\u05e0\u05e8\u05d0\u05d4 \u05de\u05db\u05d9\u05e8\u05d9\u05dd \u05d0\u05ea MCTS \u05de\u05d4\u05e4\u05e8\u05d5\u05d9\u05e7\u05d8\u05d9\u05dd \u05d4\u05de\u05e4\u05d5\u05e8\u05e1\u05de\u05d9\u05dd AlphaGo \u05d5- AlphaZero \u05e9\u05dc \u05d3\u05d9\u05e4\u05de\u05d9\u05d9\u05e0\u05d3 \u05e9\u05dc \u05d0\u05d9\u05de\u05e0\u05d5 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05de\u05e9\u05d7\u05e7 Go. \u05d0\u05e6\u05d9\u05d9\u05df AlphaZero \u05dc\u05de\u05d3 \u05dc\u05e9\u05d7\u05e7 \u05e8\u05e7 \u05d3\u05e8\u05da \u05d4\u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05e2\u05dd \u05e2\u05e6\u05de\u05d5 \u05dc\u05dc\u05d0 \u05e9\u05d5\u05dd \u05d9\u05d3\u05e2 \u05de\u05d5\u05e7\u05d3\u05dd \u05e2\u05dc Go. \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05e4\u05d5\u05ea\u05d7\u05d5 \u05d4\u05d9\u05d5 \u05db\u05d4 \u05d7\u05d6\u05e7\u05d9\u05dd \u05e9\u05d0\u05dc\u05d5\u05e3 \u05d4\u05e2\u05d5\u05dc\u05dd \u05d1-Go \u05e4\u05e8\u05e9 \u05d1\u05e2\u05e7\u05d1\u05d5\u05ea \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd (\u05dc\u05d0 \u05d6\u05d5\u05db\u05e8 \u05d0\u05d9\u05d6\u05d4). \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05de\u05ea\u05de\u05d8\u05d9 \u05de\u05d0\u05d7\u05d5\u05e8\u05d9 \u05e4\u05d9\u05ea\u05d5\u05d7 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05dc\u05d5 \u05d4\u05d9\u05d4 MCTS.\", Word.InsertLocation.replace);\nparagraphs.items[2].insertText(\"\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd MCTS \u05d4\u05d5\u05d0 \u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05de\u05e9\u05de\u05e9 \u05d1\u05e2\u05d9\u05e7\u05e8 \u05d1\u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05dc\u05e7\u05d1\u05dc\u05ea \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9\u05d5\u05ea. \u05d4\u05d5\u05d0 \u05d1\u05d5\u05e0\u05d4 \u05e2\u05e5 \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05e8\u05e6\u05ea \u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d0\u05e7\u05e8\u05d0\u05d9\u05d5\u05ea (\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d5\u05ea) \u05e8\u05d1\u05d5\u05ea \u05e9\u05dc \u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05d0\u05e4\u05e9\u05e8\u05d9\u05d9\u05dd \u05de\u05d4\u05de\u05e6\u05d1 \u05d4\u05e0\u05d5\u05db\u05d7\u05d9, \u05d5\u05de\u05e2\u05e8\u05d9\u05da \u05d0\u05ea \u05d0\u05d9\u05db\u05d5\u05ea\u05dd. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df, \u05d4\u05d5\u05d0 \u05d1\u05d5\u05d7\u05e8 \u05d0\u05ea \u05d4\u05de\u05d4\u05dc\u05da \u05e9\u05de\u05e0\u05d9\u05d1 \u05d0\u05ea \u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05d8\u05d5\u05d1\u05d5\u05ea \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d1\u05de\u05de\u05d5\u05e6\u05e2 \u05dc\u05d0\u05d5\u05e8\u05da \u05d4\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d5\u05ea. \u05d4\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05de\u05d0\u05d6\u05df \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d7\u05db\u05de\u05d4 \u05d1\u05d9\u05df \u05d7\u05e7\u05d9\u05e8\u05ea \u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05d7\u05d3\u05e9\u05d9\u05dd (exploration) \u05d4\u05e2\u05e9\u05d5\u05d9\u05d9\u05dd \u05dc\u05d4\u05ea\u05d2\u05dc\u05d5\u05ea \u05db\u05d9\u05e2\u05d9\u05dc\u05d9\u05dd, \u05dc\u05d1\u05d9\u05df \u05e0\u05d9\u05e6\u05d5\u05dc \u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05e0\u05de\u05e6\u05d0\u05d5 \u05db\u05de\u05d5\u05e6\u05dc\u05d7\u05d9\u05dd (exploitation) \u05d1\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d5\u05ea \u05e7\u05d5\u05d3\u05de\u05d5\u05ea(\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05d5\u05d1\u05d9\u05dc\u05d9\u05dd \u05dc\u05e8\u05d5\u05d1 \u05dc\u05e0\u05d9\u05e6\u05d7\u05d5\u05df \u05d1\u05de\u05e9\u05d7\u05e7).\", Word.InsertLocation.replace);\nparagraphs.items[3].insertText(\"\u05ea\u05d4\u05dc\u05d9\u05da \u05d6\u05d4 \u05d7\u05d5\u05d6\u05e8 \u05e2\u05dc \u05e2\u05e6\u05de\u05d5, \u05db\u05d0\u05e9\u05e8 \u05d1\u05db\u05dc \u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d4 \u05d4\u05e2\u05e5 \u05de\u05d5\u05e8\u05d7\u05d1 \u05d5\u05d4\u05e2\u05e8\u05db\u05d5\u05ea \u05d0\u05d9\u05db\u05d5\u05ea \u05d4\u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05de\u05ea\u05e2\u05d3\u05db\u05e0\u05d5\u05ea, \u05e2\u05d3 \u05e9\u05de\u05ea\u05e7\u05d1\u05dc\u05ea \u05d4\u05d7\u05dc\u05d8\u05d4 \u05e1\u05d5\u05e4\u05d9\u05ea. \u05d0\u05e8\u05d1\u05e2\u05ea \u05d4\u05e9\u05dc\u05d1\u05d9\u05dd \u05d4\u05de\u05e8\u05db\u05d6\u05d9\u05d9\u05dd \u05d1\u05db\u05dc \u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d4 \u05d4\u05dd: \u05d1\u05d7\u05d9\u05e8\u05ea \u05e6\u05d5\u05de\u05ea \u05d4\u05d1\u05d0 (selection), \u05d4\u05e8\u05d7\u05d1\u05ea \u05d4\u05e2\u05e5 (expansion), \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d4\u05de\u05e9\u05d7\u05e7 (simulation), \u05d5\u05e2\u05d3\u05db\u05d5\u05df \u05e2\u05e8\u05db\u05d9 \u05d4\u05e6\u05de\u05ea\u05d9\u05dd \u05e2\u05d3 \u05e9\u05d5\u05e8\u05e9 \u05d4\u05e2\u05e5 (backpropagation). \u05d4\u05e6\u05dc\u05d7\u05ea \u05d4\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05e0\u05d5\u05d1\u05e2\u05ea \u05de\u05d9\u05db\u05d5\u05dc\u05ea\u05d5 \u05dc\u05d4\u05ea\u05de\u05e7\u05d3 \u05d1\u05d0\u05d6\u05d5\u05e8\u05d9\u05dd \u05de\u05d1\u05d8\u05d9\u05d7\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e2\u05e5 \u05d4\u05d7\u05d9\u05e4\u05d5\u05e9, \u05d2\u05dd \u05d1\u05de\u05e8\u05d7\u05d1\u05d9 \u05d7\u05d9\u05e4\u05d5\u05e9 \u05e2\u05e6\u05d5\u05de\u05d9\u05dd. \u05d1\u05e1\u05d5\u05e3 \u05d4\u05de\u05d5\u05d3\u05dc, \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05de\u05d4\u05dc\u05db\u05d9 \u05de\u05e9\u05d7\u05e7 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd, (\u05de\u05e1\u05dc\u05d5\u05dc \u05d1\u05e2\u05e5) \u05d1\u05d5\u05d7\u05e8 \u05e6\u05d5\u05de\u05ea \u05d1\u05e2\u05dc \u05d4\u05e1\u05d9\u05db\u05d5\u05d9 \u05d4\u05d2\u05d1\u05d5\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc\u05e0\u05d9\u05e6\u05d7\u05d5\u05df.\", Word.InsertLocation.replace);\nparagraphs.items[4].insertText(\"\u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05d6\u05d4 \u05e7\u05e9\u05d5\u05e8 \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4. \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05d9\u05d0 \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9\u05d5\u05ea. \u05d2\u05dd \u05d1\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d0\u05e0\u05d7\u05e0\u05d5 \u05db\u05e8\u05d2\u05e2 \u05d7\u05d5\u05d6\u05d9\u05dd \u05d8\u05d5\u05e7\u05df \u05dc\u05d0\u05d7\u05e8 \u05d8\u05d5\u05e7\u05df \u05db\u05de\u05d5 \u05d1\u05de\u05e9\u05d7\u05e7 \u05d2\u05d5. \u05d1\u05e2\u05e6\u05dd \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05d2\u05d3\u05d5\u05dc \u05d1\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1-MCTS \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d4\u05d9\u05d0 \u05d1\u05e0\u05d9\u05d9\u05ea \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05d1\u05d0\u05d9\u05db\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05e4\u05d5\u05e9\u05dd \u05d1\u05e2\u05e5 \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea. \u05d0\u05d1\u05dc \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de\u05e2\u05e5 \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05e9\u05d1\u05d5 \u05d4\u05e6\u05de\u05ea\u05d9\u05dd \u05d4\u05dd \u05de\u05d4\u05dc\u05db\u05d9 \u05de\u05e9\u05d7\u05e7 \u05db\u05d0\u05df \u05db\u05dc \u05e6\u05d5\u05de\u05ea \u05d4\u05d5\u05d0 \u05e9\u05dc\u05d1 \u05d1\u05ea\u05d4\u05dc\u05d9\u05da reasoning (\u05d4\u05e0\u05de\u05e7\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc). \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d6\u05d4, \u05d1\u05e2\u05dc \u05d0\u05d9\u05db\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4\u05d4, \u05db\u05d3\u05d9 \u05dc\u05e2\u05e9\u05d5\u05ea SFT \u05dc\u05de\u05d5\u05d3\u05dc. \u05d0\u05d6 \u05d4\u05e9\u05d0\u05dc\u05d4 \u05db\u05d0\u05df \u05d0\u05d9\u05da \u05dc\u05d3\u05d2\u05d5\u05dd \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e0\u05db\u05d5\u05e0\u05d9\u05dd \u05d5\u05de\u05d2\u05d5\u05d5\u05e0\u05d9\u05dd \u05e2\u05dd \u05d2\u05d9\u05e9\u05d4 \u05d6\u05d5?\", Word.InsertLocation.replace);\nparagraphs.items[5].insertText(\"\u05db\u05d0\u05de\u05d5\u05e8 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea MCTS \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05e8\u05d5\u05e8\u05d4 (\u05d4\u05d0\u05dd \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05e0\u05db\u05d5\u05df \u05d0\u05d5 \u05dc\u05d0) \u05d1\u05e1\u05d5\u05e3 \u05d4\u05d2\u05e0\u05e8\u05d5\u05d8. \u05dc\u05e2\u05d5\u05de\u05ea \u05d6\u05d0\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc (reward) \u05d1\u05d0\u05de\u05e6\u05e2 \u05e9\u05e8\u05e9\u05e8\u05ea \u05d4\u05d4\u05e0\u05de\u05e7\u05d4 \u05d4\u05d5\u05d0 \u05de\u05e9\u05d4\u05d5 \u05d1\u05e8\u05d5\u05e8 (\u05d3\u05f4\u05d0 \u05d9\u05e9 \u05d1-PPO \u05d0\u05ea \u05d0\u05d5\u05ea\u05d4 \u05d4\u05d1\u05e2\u05d9\u05d4 - \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea reward \u05e9\u05d0\u05d9\u05de\u05e0\u05d5 \u05d0\u05d5\u05dc\u05dd \u05d4\u05d9\u05d0 \u05e0\u05d5\u05ea\u05e0\u05ea \u05e6\u05d9\u05d5\u05df \u05dc\u05db\u05dc \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05d5\u05dc\u05d0 \u05dc\u05d7\u05dc\u05e7\u05d5 \u05d5\u05d0\u05d6 \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea value \u05d4\u05de\u05e9\u05e2\u05e8\u05db\u05ea \u05d0\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05e9\u05dc\u05d1\u05d9 \u05d1\u05d9\u05e0\u05d9\u05d9\u05dd - \u05d3\u05e8\u05da \u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05e2\u05d9\u05d9\u05ea \u05e8\u05d2\u05e8\u05e1\u05d9\u05d4). \u05d1-MCTS \u05d1\u05e0\u05d9\u05d9\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05d4\u05de\u05e7\u05e0\u05d4 \u05e6\u05d9\u05d5\u05df \u05dc\u05e6\u05d5\u05de\u05ea (\u05e4\u05ea\u05e8\u05d5\u05df \u05d7\u05dc\u05e7\u05d9) \u05d4\u05d5\u05d0 \u05e7\u05e8\u05d9\u05d8\u05d9 \u05db\u05d9 \u05d0\u05d7\u05e8\u05ea \u05dc\u05d0 \u05e0\u05e6\u05dc\u05d9\u05d7 \u05dc\u05d1\u05e0\u05d5\u05ea \u05d0\u05ea \u05e2\u05e5 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d8\u05d5\u05d1\u05d4 (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e0\u05d9\u05d1\u05d4 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d8\u05d5\u05d1\u05d9\u05dd \u05dc\u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea). \u05db\u05de\u05d5\u05d1\u05df \u05db\u05dc \u05e6\u05d5\u05de\u05ea \u05d1\u05e2\u05e5 \u05e0\u05d1\u05e0\u05d4 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4.\", Word.InsertLocation.replace);\nparagraphs.items[6].insertText(\"\u05d1\u05d4\u05ea\u05d7\u05dc\u05d4 \u05e6\u05d9\u05d5\u05df \u05d4\u05e6\u05d5\u05de\u05ea (= \u05e4\u05ea\u05e8\u05d5\u05df \u05d7\u05dc\u05e7\u05d9 \u05e2\u05d3 \u05e9\u05dc\u05d1 \u05de\u05e1\u05d5\u05d9\u05dd) \u05e0\u05d1\u05e0\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e9\u05db\u05d9\u05d7\u05d5\u05ea \u05d4\u05d5\u05e4\u05e2\u05ea\u05d5 \u05d1\u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e0\u05db\u05d5\u05e0\u05d9\u05dd \u05e9\u05dc \u05d4\u05d1\u05e2\u05d9\u05d4. \u05db\u05db\u05dc \u05d4\u05d5\u05d0 \u05de\u05d5\u05e4\u05d9\u05e2 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e9\u05e8\u05e9\u05e8\u05d0\u05d5\u05ea \u05d4\u05e0\u05de\u05e7\u05d4 \u05d4\u05de\u05d5\u05d1\u05d9\u05dc\u05d5\u05ea \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05e0\u05db\u05d5\u05df, \u05e6\u05d9\u05d5\u05df \u05e9\u05dc\u05d5 \u05d2\u05d1\u05d5\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05d1\u05e9\u05dc\u05d1\u05d9\u05dd \u05de\u05d0\u05d5\u05d7\u05e8 \u05d9\u05d5\u05ea\u05e8 (\u05db\u05d0\u05e9\u05e8 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05e6\u05d9\u05d5\u05df \u05de\u05e6\u05d9\u05d9\u05e6\u05d1\u05ea) \u05d4\u05de\u05d0\u05de\u05e8 \u05e2\u05d5\u05e9\u05d9\u05dd \u05de\u05e9\u05d4\u05d5 \u05d3\u05d5\u05de\u05d4 \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc \u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05d0\u05d9\u05de\u05d5\u05df RLHF \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d1\u05db\u05dc \u05e2\u05d5\u05de\u05e7 (\u05e9\u05db\u05d1\u05d4) \u05e9\u05dc \u05e2\u05e5 \u05dc\u05d5\u05e7\u05d7\u05d9\u05dd \u05e6\u05de\u05ea\u05d9\u05dd \u05d1\u05e2\u05dc\u05d9 \u05e6\u05d9\u05d5\u05e0\u05d9\u05dd \u05d4\u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d5\u05d4\u05e0\u05de\u05d5\u05db\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d5\u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e6\u05d9\u05d5\u05df \u05e6\u05d5\u05de\u05ea \u05d1\u05e1\u05d2\u05e0\u05d5\u05df Bradley-Terry (\u05db\u05de\u05d5 \u05e9\u05de\u05e7\u05d5\u05d1\u05dc \u05d1-RLHF \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9). \u05db\u05d0\u05de\u05d5\u05e8 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05e6\u05d9\u05d5\u05df \u05de\u05e9\u05de\u05e9\u05ea \u05d0\u05d5\u05ea\u05e0\u05d5 \u05dc\u05d1\u05d7\u05d9\u05e8\u05d4 \u05de\u05d0\u05d9\u05d6\u05d4 \u05e6\u05d5\u05de\u05ea \u05dc\u05d3\u05d2\u05d5\u05dd \u05e9\u05dc\u05d1 \u05d4\u05d1\u05d0 \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9 (Upper Confidence bounds for Trees (UCT \u05d4\u05de\u05d2'\u05e0\u05d2\u05dc \u05d1\u05d9\u05df exploration vs exploitation.\", Word.InsertLocation.replace);\nparagraphs.items[7].insertText(\"\u05db\u05d3\u05d9 \u05dc\u05d4\u05d2\u05d9\u05e2 \u05dc\u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d0\u05d9\u05db\u05d5\u05ea\u05d9\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d4\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05ea\u05d1\u05e7\u05e9 \u05dc\u05de\u05de\u05e9 \u05db\u05dc \u05e9\u05dc\u05d1 \u05d1\u05e9\u05e8\u05e9\u05e8\u05ea \u05d4\u05e0\u05de\u05e7\u05d4 \u05d1\u05e4\u05d9\u05d9\u05d8\u05d5\u05df \u05d5\u05d0\u05dd \u05e7\u05d5\u05d3 \u05d6\u05d4 \u05dc\u05d0 \u05e2\u05d5\u05d1\u05e8 \u05d8\u05e1\u05d8\u05d9\u05dd, \u05d4\u05e6\u05d5\u05de\u05ea \u05e0\u05e4\u05e1\u05dc. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05ea\u05d7\u05d9\u05dc \u05de\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05e7\u05d8\u05df, \u05d9\u05d5\u05e6\u05e8 \u05e2\u05e5 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea (\u05e2\u05dd \u05db\u05dc \u05d4\u05e9\u05dc\u05d1\u05d9\u05dd \u05e9\u05ea\u05d9\u05d0\u05e8\u05ea\u05d9), \u05d1\u05d5\u05d7\u05e8 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d4\u05db\u05d9 \u05d0\u05d9\u05db\u05d5\u05ea\u05d9\u05d9\u05dd (\u05d1\u05e2\u05dc\u05d9 \u05e6\u05d9\u05d5\u05e0\u05d9 \u05d4\u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8), \u05de\u05e6\u05d1\u05e2 SFT \u05e2\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d5\u05d7\u05d5\u05d6\u05e8 \u05e2\u05dc \u05d6\u05d4 \u05e2\u05d5\u05d3 \u05e4\u05e2\u05dd. \u05d5\u05db\u05ea\u05d5\u05e6\u05d0\u05d4 \u05de\u05db\u05da \u05d0\u05e0\u05d5 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e7\u05d8\u05df \u05d5\u05d7\u05de\u05d5\u05d3 \u05d0\u05d1\u05dc \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea \u05d3\u05d9 \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea (\u05dc\u05db\u05d0\u05d5\u05e8\u05d4).\", Word.InsertLocation.replace);\nparagraphs.items[10].insertText(\"https://arxiv.org/abs/2501.04519\", Word.InsertLocation.replace);\n\n// Delete obsolete paragraphs (in reverse order to keep indices stable)\nparagraphs.items[9].delete();\nparagraphs.items[8].delete();\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\nif ($d.Paragraphs.Count -ne 11) {\n    throw \"Unexpected paragraph count: $($d.Paragraphs.Count)\"\n}\n\n# Update paragraph text in place (1-based COM indices; stable while only\n# mutating Range.Text without inserting/removing paragraph marks).\n$d.Paragraphs.Item(1).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7: 24.05.25\" + [char]0x0B + \"rStar-Math: Small LLMs Can Master Math Reasoning with Self-Evolved Deep Thinking\"\n$d.Paragraphs.Item(2).Range.Text = \"\u05db\u05de\u05d4 \u05d9\u05de\u05d9\u05dd \u05dc\u05d0 \u05e1\u05e7\u05e8\u05ea\u05d9 \u05de\u05d0\u05de\u05e8 \u05d0\u05d1\u05dc \u05d1\u05d9\u05d5\u05dd \u05d4\u05d5\u05dc\u05d3\u05ea\u05d9 \u05dc\u05d0 \u05d9\u05db\u05d5\u05dc\u05ea\u05d9 \u05dc\u05d0 \u05dc\u05db\u05ea\u05d5\u05d1 \u05e1\u05e7\u05d9\u05e8\u05d4 \u05dc\u05de\u05e8\u05d5\u05ea \u05d4\u05e2\u05d5\u05de\u05e1 \u05d4\u05de\u05d8\u05d5\u05e8\u05e3. \u05d4\u05d9\u05d5\u05dd \u05d0\u05e1\u05e7\u05d5\u05e8 \u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e9\u05d9\u05e6\u05d0 \u05dc\u05e4\u05e0\u05d9 4 \u05d7\u05d5\u05d3\u05e9\u05d9\u05dd \u05d5\u05d4\u05d5\u05d0 \u05de\u05e9\u05dc\u05d1 \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05dc\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea MCTS \u05e9\u05d6\u05d4 \u05e7\u05d9\u05e6\u05d5\u05e8 \u05e9\u05dc Monte Carlo Tree Search. \u05e8\u05d5\u05d1\u05db\u05dd \u05db\u05e0\u05e8\u05d0\u05d4 \u05de\u05db\u05d9\u05e8\u05d9\u05dd \u05d0\u05ea MCTS \u05de\u05d4\u05e4\u05e8\u05d5\u05d9\u05e7\u05d8\u05d9\u05dd \u05d4\u05de\u05e4\u05d5\u05e8\u05e1\u05de\u05d9\u05dd AlphaGo \u05d5- AlphaZero \u05e9\u05dc \u05d3\u05d9\u05e4\u05de\u05d9\u05d9\u05e0\u05d3 \u05e9\u05dc \u05d0\u05d9\u05de\u05e0\u05d5 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05de\u05e9\u05d7\u05e7 Go. \u05d0\u05e6\u05d9\u05d9\u05df AlphaZero \u05dc\u05de\u05d3 \u05dc\u05e9\u05d7\u05e7 \u05e8\u05e7 \u05d3\u05e8\u05da \u05d4\u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05e2\u05dd \u05e2\u05e6\u05de\u05d5 \u05dc\u05dc\u05d0 \u05e9\u05d5\u05dd \u05d9\u05d3\u05e2 \u05de\u05d5\u05e7\u05d3\u05dd \u05e2\u05dc Go. \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05e4\u05d5\u05ea\u05d7\u05d5 \u05d4\u05d9\u05d5 \u05db\u05d4 \u05d7\u05d6\u05e7\u05d9\u05dd \u05e9\u05d0\u05dc\u05d5\u05e3 \u05d4\u05e2\u05d5\u05dc\u05dd \u05d1-Go \u05e4\u05e8\u05e9 \u05d1\u05e2\u05e7\u05d1\u05d5\u05ea \u05d0\u05d7\u05d3 \u05de\u05d4\u05dd (\u05dc\u05d0 \u05d6\u05d5\u05db\u05e8 \u05d0\u05d9\u05d6\u05d4). \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05de\u05ea\u05de\u05d8\u05d9 \u05de\u05d0\u05d7\u05d5\u05e8\u05d9 \u05e4\u05d9\u05ea\u05d5\u05d7 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05dc\u05d5 \u05d4\u05d9\u05d4 MCTS.\"\n$d.Paragraphs.Item(3).Range.Text = \"\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd MCTS \u05d4\u05d5\u05d0 \u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05de\u05e9\u05de\u05e9 \u05d1\u05e2\u05d9\u05e7\u05e8 \u05d1\u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05dc\u05e7\u05d1\u05dc\u05ea \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9\u05d5\u05ea. \u05d4\u05d5\u05d0 \u05d1\u05d5\u05e0\u05d4 \u05e2\u05e5 \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05e8\u05e6\u05ea \u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d0\u05e7\u05e8\u05d0\u05d9\u05d5\u05ea (\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d5\u05ea) \u05e8\u05d1\u05d5\u05ea \u05e9\u05dc \u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05d0\u05e4\u05e9\u05e8\u05d9\u05d9\u05dd \u05de\u05d4\u05de\u05e6\u05d1 \u05d4\u05e0\u05d5\u05db\u05d7\u05d9, \u05d5\u05de\u05e2\u05e8\u05d9\u05da \u05d0\u05ea \u05d0\u05d9\u05db\u05d5\u05ea\u05dd. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df, \u05d4\u05d5\u05d0 \u05d1\u05d5\u05d7\u05e8 \u05d0\u05ea \u05d4\u05de\u05d4\u05dc\u05da \u05e9\u05de\u05e0\u05d9\u05d1 \u05d0\u05ea \u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d4\u05d8\u05d5\u05d1\u05d5\u05ea \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d1\u05de\u05de\u05d5\u05e6\u05e2 \u05dc\u05d0\u05d5\u05e8\u05da \u05d4\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d5\u05ea. \u05d4\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05de\u05d0\u05d6\u05df \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d7\u05db\u05de\u05d4 \u05d1\u05d9\u05df \u05d7\u05e7\u05d9\u05e8\u05ea \u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05d7\u05d3\u05e9\u05d9\u05dd (exploration) \u05d4\u05e2\u05e9\u05d5\u05d9\u05d9\u05dd \u05dc\u05d4\u05ea\u05d2\u05dc\u05d5\u05ea \u05db\u05d9\u05e2\u05d9\u05dc\u05d9\u05dd, \u05dc\u05d1\u05d9\u05df \u05e0\u05d9\u05e6\u05d5\u05dc \u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05e0\u05de\u05e6\u05d0\u05d5 \u05db\u05de\u05d5\u05e6\u05dc\u05d7\u05d9\u05dd (exploitation) \u05d1\u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d5\u05ea \u05e7\u05d5\u05d3\u05de\u05d5\u05ea(\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05d5\u05d1\u05d9\u05dc\u05d9\u05dd \u05dc\u05e8\u05d5\u05d1 \u05dc\u05e0\u05d9\u05e6\u05d7\u05d5\u05df \u05d1\u05de\u05e9\u05d7\u05e7).\"\n$d.Paragraphs.Item(4).Range.Text = \"\u05ea\u05d4\u05dc\u05d9\u05da \u05d6\u05d4 \u05d7\u05d5\u05d6\u05e8 \u05e2\u05dc \u05e2\u05e6\u05de\u05d5, \u05db\u05d0\u05e9\u05e8 \u05d1\u05db\u05dc \u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d4 \u05d4\u05e2\u05e5 \u05de\u05d5\u05e8\u05d7\u05d1 \u05d5\u05d4\u05e2\u05e8\u05db\u05d5\u05ea \u05d0\u05d9\u05db\u05d5\u05ea \u05d4\u05de\u05d4\u05dc\u05db\u05d9\u05dd \u05de\u05ea\u05e2\u05d3\u05db\u05e0\u05d5\u05ea, \u05e2\u05d3 \u05e9\u05de\u05ea\u05e7\u05d1\u05dc\u05ea \u05d4\u05d7\u05dc\u05d8\u05d4 \u05e1\u05d5\u05e4\u05d9\u05ea. \u05d0\u05e8\u05d1\u05e2\u05ea \u05d4\u05e9\u05dc\u05d1\u05d9\u05dd \u05d4\u05de\u05e8\u05db\u05d6\u05d9\u05d9\u05dd \u05d1\u05db\u05dc \u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d4 \u05d4\u05dd: \u05d1\u05d7\u05d9\u05e8\u05ea \u05e6\u05d5\u05de\u05ea \u05d4\u05d1\u05d0 (selection), \u05d4\u05e8\u05d7\u05d1\u05ea \u05d4\u05e2\u05e5 (expansion), \u05e1\u05d9\u05de\u05d5\u05dc\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d4\u05de\u05e9\u05d7\u05e7 (simulation), \u05d5\u05e2\u05d3\u05db\u05d5\u05df \u05e2\u05e8\u05db\u05d9 \u05d4\u05e6\u05de\u05ea\u05d9\u05dd \u05e2\u05d3 \u05e9\u05d5\u05e8\u05e9 \u05d4\u05e2\u05e5 (backpropagation). \u05d4\u05e6\u05dc\u05d7\u05ea \u05d4\u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05e0\u05d5\u05d1\u05e2\u05ea \u05de\u05d9\u05db\u05d5\u05dc\u05ea\u05d5 \u05dc\u05d4\u05ea\u05de\u05e7\u05d3 \u05d1\u05d0\u05d6\u05d5\u05e8\u05d9\u05dd \u05de\u05d1\u05d8\u05d9\u05d7\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e2\u05e5 \u05d4\u05d7\u05d9\u05e4\u05d5\u05e9, \u05d2\u05dd \u05d1\u05de\u05e8\u05d7\u05d1\u05d9 \u05d7\u05d9\u05e4\u05d5\u05e9 \u05e2\u05e6\u05d5\u05de\u05d9\u05dd. \u05d1\u05e1\u05d5\u05e3 \u05d4\u05de\u05d5\u05d3\u05dc, \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05de\u05d4\u05dc\u05db\u05d9 \u05de\u05e9\u05d7\u05e7 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd, (\u05de\u05e1\u05dc\u05d5\u05dc \u05d1\u05e2\u05e5) \u05d1\u05d5\u05d7\u05e8 \u05e6\u05d5\u05de\u05ea \u05d1\u05e2\u05dc \u05d4\u05e1\u05d9\u05db\u05d5\u05d9 \u05d4\u05d2\u05d1\u05d5\u05d4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc\u05e0\u05d9\u05e6\u05d7\u05d5\u05df.\"\n$d.Paragraphs.Item(5).Range.Text = \"\u05d0\u05d1\u05dc \u05d0\u05d9\u05da \u05d6\u05d4 \u05e7\u05e9\u05d5\u05e8 \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4. \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05d9\u05d0 \u05d0\u05d5\u05d8\u05d5\u05e8\u05d2\u05e8\u05e1\u05d9\u05d1\u05d9\u05d5\u05ea. \u05d2\u05dd \u05d1\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d0\u05e0\u05d7\u05e0\u05d5 \u05db\u05e8\u05d2\u05e2 \u05d7\u05d5\u05d6\u05d9\u05dd \u05d8\u05d5\u05e7\u05df \u05dc\u05d0\u05d7\u05e8 \u05d8\u05d5\u05e7\u05df \u05db\u05de\u05d5 \u05d1\u05de\u05e9\u05d7\u05e7 \u05d2\u05d5. \u05d1\u05e2\u05e6\u05dd \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05d4\u05d2\u05d3\u05d5\u05dc \u05d1\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1-MCTS \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d4\u05d9\u05d0 \u05d1\u05e0\u05d9\u05d9\u05ea \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05d1\u05d0\u05d9\u05db\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05e4\u05d5\u05e9\u05dd \u05d1\u05e2\u05e5 \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea. \u05d0\u05d1\u05dc \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de\u05e2\u05e5 \u05d4\u05d7\u05dc\u05d8\u05d5\u05ea \u05e9\u05d1\u05d5 \u05d4\u05e6\u05de\u05ea\u05d9\u05dd \u05d4\u05dd \u05de\u05d4\u05dc\u05db\u05d9 \u05de\u05e9\u05d7\u05e7 \u05db\u05d0\u05df \u05db\u05dc \u05e6\u05d5\u05de\u05ea \u05d4\u05d5\u05d0 \u05e9\u05dc\u05d1 \u05d1\u05ea\u05d4\u05dc\u05d9\u05da reasoning (\u05d4\u05e0\u05de\u05e7\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc). \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d6\u05d4, \u05d1\u05e2\u05dc \u05d0\u05d9\u05db\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4\u05d4, \u05db\u05d3\u05d9 \u05dc\u05e2\u05e9\u05d5\u05ea SFT \u05dc\u05de\u05d5\u05d3\u05dc. \u05d0\u05d6 \u05d4\u05e9\u05d0\u05dc\u05d4 \u05db\u05d0\u05df \u05d0\u05d9\u05da \u05dc\u05d3\u05d2\u05d5\u05dd \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e0\u05db\u05d5\u05e0\u05d9\u05dd \u05d5\u05de\u05d2\u05d5\u05d5\u05e0\u05d9\u05dd \u05e2\u05dd \u05d2\u05d9\u05e9\u05d4 \u05d6\u05d5?\"\n$d.Paragraphs.Item(6).Range.Text = \"\u05db\u05d0\u05de\u05d5\u05e8 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea MCTS \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05e8\u05d5\u05e8\u05d4 (\u05d4\u05d0\u05dd \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05e0\u05db\u05d5\u05df \u05d0\u05d5 \u05dc\u05d0) \u05d1\u05e1\u05d5\u05e3 \u05d4\u05d2\u05e0\u05e8\u05d5\u05d8. \u05dc\u05e2\u05d5\u05de\u05ea \u05d6\u05d0\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc (reward) \u05d1\u05d0\u05de\u05e6\u05e2 \u05e9\u05e8\u05e9\u05e8\u05ea \u05d4\u05d4\u05e0\u05de\u05e7\u05d4 \u05d4\u05d5\u05d0 \u05de\u05e9\u05d4\u05d5 \u05d1\u05e8\u05d5\u05e8 (\u05d3\u05f4\u05d0 \u05d9\u05e9 \u05d1-PPO \u05d0\u05ea \u05d0\u05d5\u05ea\u05d4 \u05d4\u05d1\u05e2\u05d9\u05d4 - \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea reward \u05e9\u05d0\u05d9\u05de\u05e0\u05d5 \u05d0\u05d5\u05dc\u05dd \u05d4\u05d9\u05d0 \u05e0\u05d5\u05ea\u05e0\u05ea \u05e6\u05d9\u05d5\u05df \u05dc\u05db\u05dc \u05d4\u05e4\u05ea\u05e8\u05d5\u05df \u05d5\u05dc\u05d0 \u05dc\u05d7\u05dc\u05e7\u05d5 \u05d5\u05d0\u05d6 \u05d0\u05e0\u05d5 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea value \u05d4\u05de\u05e9\u05e2\u05e8\u05db\u05ea \u05d0\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05e9\u05dc\u05d1\u05d9 \u05d1\u05d9\u05e0\u05d9\u05d9\u05dd - \u05d3\u05e8\u05da \u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05e2\u05d9\u05d9\u05ea \u05e8\u05d2\u05e8\u05e1\u05d9\u05d4). \u05d1-MCTS \u05d1\u05e0\u05d9\u05d9\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05d4\u05de\u05e7\u05e0\u05d4 \u05e6\u05d9\u05d5\u05df \u05dc\u05e6\u05d5\u05de\u05ea (\u05e4\u05ea\u05e8\u05d5\u05df \u05d7\u05dc\u05e7\u05d9) \u05d4\u05d5\u05d0 \u05e7\u05e8\u05d9\u05d8\u05d9 \u05db\u05d9 \u05d0\u05d7\u05e8\u05ea \u05dc\u05d0 \u05e0\u05e6\u05dc\u05d9\u05d7 \u05dc\u05d1\u05e0\u05d5\u05ea \u05d0\u05ea \u05e2\u05e5 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d8\u05d5\u05d1\u05d4 (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e0\u05d9\u05d1\u05d4 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d8\u05d5\u05d1\u05d9\u05dd \u05dc\u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea). \u05db\u05de\u05d5\u05d1\u05df \u05db\u05dc \u05e6\u05d5\u05de\u05ea \u05d1\u05e2\u05e5 \u05e0\u05d1\u05e0\u05d4 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d3\u05d2\u05d9\u05de\u05d4 \u05de\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4.\"\n$d.Paragraphs.Item(7).Range.Text = \"\u05d1\u05d4\u05ea\u05d7\u05dc\u05d4 \u05e6\u05d9\u05d5\u05df \u05d4\u05e6\u05d5\u05de\u05ea (= \u05e4\u05ea\u05e8\u05d5\u05df \u05d7\u05dc\u05e7\u05d9 \u05e2\u05d3 \u05e9\u05dc\u05d1 \u05de\u05e1\u05d5\u05d9\u05dd) \u05e0\u05d1\u05e0\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e9\u05db\u05d9\u05d7\u05d5\u05ea \u05d4\u05d5\u05e4\u05e2\u05ea\u05d5 \u05d1\u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e0\u05db\u05d5\u05e0\u05d9\u05dd \u05e9\u05dc \u05d4\u05d1\u05e2\u05d9\u05d4. \u05db\u05db\u05dc \u05d4\u05d5\u05d0 \u05de\u05d5\u05e4\u05d9\u05e2 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e9\u05e8\u05e9\u05e8\u05d0\u05d5\u05ea \u05d4\u05e0\u05de\u05e7\u05d4 \u05d4\u05de\u05d5\u05d1\u05d9\u05dc\u05d5\u05ea \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05e0\u05db\u05d5\u05df, \u05e6\u05d9\u05d5\u05df \u05e9\u05dc\u05d5 \u05d2\u05d1\u05d5\u05d4 \u05d9\u05d5\u05ea\u05e8. \u05d1\u05e9\u05dc\u05d1\u05d9\u05dd \u05de\u05d0\u05d5\u05d7\u05e8 \u05d9\u05d5\u05ea\u05e8 (\u05db\u05d0\u05e9\u05e8 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05e6\u05d9\u05d5\u05df \u05de\u05e6\u05d9\u05d9\u05e6\u05d1\u05ea) \u05d4\u05de\u05d0\u05de\u05e8 \u05e2\u05d5\u05e9\u05d9\u05dd \u05de\u05e9\u05d4\u05d5 \u05d3\u05d5\u05de\u05d4 \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc \u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05d0\u05d9\u05de\u05d5\u05df RLHF \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d1\u05db\u05dc \u05e2\u05d5\u05de\u05e7 (\u05e9\u05db\u05d1\u05d4) \u05e9\u05dc \u05e2\u05e5 \u05dc\u05d5\u05e7\u05d7\u05d9\u05dd \u05e6\u05de\u05ea\u05d9\u05dd \u05d1\u05e2\u05dc\u05d9 \u05e6\u05d9\u05d5\u05e0\u05d9\u05dd \u05d4\u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d5\u05d4\u05e0\u05de\u05d5\u05db\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d5\u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e6\u05d9\u05d5\u05df \u05e6\u05d5\u05de\u05ea \u05d1\u05e1\u05d2\u05e0\u05d5\u05df Bradley-Terry (\u05db\u05de\u05d5 \u05e9\u05de\u05e7\u05d5\u05d1\u05dc \u05d1-RLHF \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9). \u05db\u05d0\u05de\u05d5\u05e8 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05e6\u05d9\u05d5\u05df \u05de\u05e9\u05de\u05e9\u05ea \u05d0\u05d5\u05ea\u05e0\u05d5 \u05dc\u05d1\u05d7\u05d9\u05e8\u05d4 \u05de\u05d0\u05d9\u05d6\u05d4 \u05e6\u05d5\u05de\u05ea \u05dc\u05d3\u05d2\u05d5\u05dd \u05e9\u05dc\u05d1 \u05d4\u05d1\u05d0 \u05dc\u05e4\u05ea\u05e8\u05d5\u05df \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d0\u05dc\u05d2\u05d5\u05e8\u05d9\u05ea\u05dd \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9 (Upper Confidence bounds for Trees (UCT \u05d4\u05de\u05d2'\u05e0\u05d2\u05dc \u05d1\u05d9\u05df exploration vs exploitation.\"\n$d.Paragraphs.Item(8).Range.Text = \"\u05db\u05d3\u05d9 \u05dc\u05d4\u05d2\u05d9\u05e2 \u05dc\u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d0\u05d9\u05db\u05d5\u05ea\u05d9\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d4\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05ea\u05d1\u05e7\u05e9 \u05dc\u05de\u05de\u05e9 \u05db\u05dc \u05e9\u05dc\u05d1 \u05d1\u05e9\u05e8\u05e9\u05e8\u05ea \u05d4\u05e0\u05de\u05e7\u05d4 \u05d1\u05e4\u05d9\u05d9\u05d8\u05d5\u05df \u05d5\u05d0\u05dd \u05e7\u05d5\u05d3 \u05d6\u05d4 \u05dc\u05d0 \u05e2\u05d5\u05d1\u05e8 \u05d8\u05e1\u05d8\u05d9\u05dd, \u05d4\u05e6\u05d5\u05de\u05ea \u05e0\u05e4\u05e1\u05dc. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05ea\u05d7\u05d9\u05dc \u05de\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05e7\u05d8\u05df, \u05d9\u05d5\u05e6\u05e8 \u05e2\u05e5 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea (\u05e2\u05dd \u05db\u05dc \u05d4\u05e9\u05dc\u05d1\u05d9\u05dd \u05e9\u05ea\u05d9\u05d0\u05e8\u05ea\u05d9), \u05d1\u05d5\u05d7\u05e8 \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05d4\u05db\u05d9 \u05d0\u05d9\u05db\u05d5\u05ea\u05d9\u05d9\u05dd (\u05d1\u05e2\u05dc\u05d9 \u05e6\u05d9\u05d5\u05e0\u05d9 \u05d4\u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8), \u05de\u05e6\u05d1\u05e2 SFT \u05e2\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05d5\u05d7\u05d5\u05d6\u05e8 \u05e2\u05dc \u05d6\u05d4 \u05e2\u05d5\u05d3 \u05e4\u05e2\u05dd. \u05d5\u05db\u05ea\u05d5\u05e6\u05d0\u05d4 \u05de\u05db\u05da \u05d0\u05e0\u05d5 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e7\u05d8\u05df \u05d5\u05d7\u05de\u05d5\u05d3 \u05d0\u05d1\u05dc \u05de\u05e1\u05d5\u05d2\u05dc \u05dc\u05e4\u05ea\u05d5\u05e8 \u05d1\u05e2\u05d9\u05d5\u05ea \u05de\u05ea\u05de\u05d8\u05d9\u05d5\u05ea \u05d3\u05d9 \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea (\u05dc\u05db\u05d0\u05d5\u05e8\u05d4).\"\n$d.Paragraphs.Item(11).Range.Text = \"https://arxiv.org/abs/2501.04519\"\n\n# Remove the two obsolete paragraphs, highest index first so earlier\n# indices used above stay valid.\n$d.Paragraphs.Item(10).Range.Delete()\n$d.Paragraphs.Item(9).Range.Delete()\n"}
